# Add 3 new city rows (170-172) to Sheet1, matching the upstream "Add files
# via upload" commit: Yadamah, Al Artawiah, AL Bejadiah.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (borders/style) of the last existing data row (169)
# down onto the three new rows so the new cells pick up the same style
# index as the rest of the table.
$ws.Range("A169:G169").Copy()
$ws.Range("A170:G172").PasteSpecial(-4122)

# Row 170: Yadamah / يدمة - Najran region / South of the Kingdom
$ws.Cells.Item(170, 1).Value = "Yadamah"
$ws.Cells.Item(170, 2).Value = "Yadamah"
$ws.Cells.Item(170, 3).Value = "يدمة"
$ws.Cells.Item(170, 4).Value = 18.532356
$ws.Cells.Item(170, 5).Value = 44.228935
$ws.Cells.Item(170, 6).Value = "منطقة نجران"
$ws.Cells.Item(170, 7).Value = "جنوب المملكة"

# Row 171: Al Artawiah / الأرطاية - Riyadh region / Center of the Kingdom
$ws.Cells.Item(171, 1).Value = "Al Artawiah"
$ws.Cells.Item(171, 2).Value = "Al Artawiah"
$ws.Cells.Item(171, 3).Value = "الأرطاية"
$ws.Cells.Item(171, 4).Value = 26.50162
$ws.Cells.Item(171, 5).Value = 45.345532
$ws.Cells.Item(171, 6).Value = "منطقة الرياض"
$ws.Cells.Item(171, 7).Value = "وسط المملكة"

# Row 172: AL Bejadiah / البجادية - Riyadh region / Center of the Kingdom
$ws.Cells.Item(172, 1).Value = "AL Bejadiah"
$ws.Cells.Item(172, 2).Value = "AL Bejadiah"
$ws.Cells.Item(172, 3).Value = "البجادية"
$ws.Cells.Item(172, 4).Value = 23.61806
$ws.Cells.Item(172, 5).Value = 45.388869
$ws.Cells.Item(172, 6).Value = "منطقة الرياض"
$ws.Cells.Item(172, 7).Value = "وسط المملكة"

# Extend the sheet selection to cover the new rows (A1:G169 -> A1:G172),
# matching the saved <selection sqref="..."/> in the workbook view.
$ws.Range("A1:G172").Select()
